# Insert a new daily price record at the top of the data block (row 9),
# shifting all existing records down by one row, and populate the new
# row with the latest report's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 9 - this shifts rows 9:97 down
# to 10:98 and extends the used range to A1:T98.
$ws.Rows("9:9").Insert()

# Populate the newly inserted row 9 with the new record.
$ws.Range("A9").Value = 1
$ws.Range("B9").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C9").Value = "Arica y Parinacota"
$ws.Range("D9").Value = 44685
$ws.Range("E9").Value = 15
$ws.Range("F9").Value = "Fruta"
$ws.Range("G9").Value = 100102
$ws.Range("H9").Value = "Cítricos"
$ws.Range("I9").Value = 100102004
$ws.Range("J9").Value = "Mandarina"
$ws.Range("K9").Value = "Murcott"
$ws.Range("L9").Value = "Segunda"
$ws.Range("M9").Value = 300
$ws.Range("N9").Value = 18000
$ws.Range("O9").Value = 19000
$ws.Range("P9").Value = 18500
$ws.Range("Q9").Value = "$/caja 20 kilos"
$ws.Range("R9").Value = "Región de Coquimbo"
$ws.Range("S9").Value = 925
$ws.Range("T9").Value = 20
